# ---------------------------------------------------------------------------
# ISYS3001.docx edit:
#   1) "A cool quote by Edsger Dijkstra:" -> wrap "Edsger" in its own run
#      flanked by proofErr spellStart/spellEnd markers.
#   2) Quote paragraph -> move the Arial/italic run formatting onto the
#      paragraph mark (pPr/rPr) and merge the closing curly quote into the
#      same run as the quote text.
#   3) Append a new "Learning version management is important" paragraph
#      (Arial/21/shading, no italics) after the quote, followed by a
#      trailing empty paragraph.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function Get-ParagraphIndexContainingText($doc, $text) {
    $rng = $doc.Content
    $ok = $rng.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) { throw "Get-ParagraphIndexContainingText: text not found: $text" }
    return $rng.Paragraphs(1).Index
}

# --- 1) split the "Edsger Dijkstra" quote-attribution paragraph -----------
$idxQuoteBy = Get-ParagraphIndexContainingText $d "A cool quote by Edsger Dijkstra:"
$pQuoteBy = $d.Paragraphs($idxQuoteBy).Range

$xmlQuoteBy = @'
<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">A cool quote by </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Edsger</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Dijkstra:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$pQuoteBy.InsertXML($xmlQuoteBy)

# --- 2) re-shape the quote paragraph ---------------------------------------
$idxQuote = Get-ParagraphIndexContainingText $d "Computer science is no more about computers"
$pQuote = $d.Paragraphs($idxQuote).Range

$xmlQuote = @'
<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:iCs/><w:color w:val="4D5156"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:t>“</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:iCs/><w:color w:val="4D5156"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Computer science is no more about computers than astronomy is about telescopes.”</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$pQuote.InsertXML($xmlQuote)

# --- 3) append the "Learning version management" paragraph + trailing blank
$idxQuote = Get-ParagraphIndexContainingText $d "Computer science is no more about computers"
$pQuote = $d.Paragraphs($idxQuote).Range
$null = $pQuote.InsertParagraphAfter()

$pNew = $d.Paragraphs($idxQuote + 1).Range

$xmlLearning = @'
<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Learning version management is important</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$pNew.InsertXML($xmlLearning)
